$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 57
$lastRow = $newRow - 1

# Copy the formatting of the last existing row onto the new row so the
# appended log entry matches the sheet's existing style (s="3": centered,
# default font/border/fill).
$ws.Range("A" + $lastRow + ":H" + $lastRow).Copy()
$ws.Range("A" + $newRow + ":H" + $newRow).PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Cells.Item($newRow, 1).Value = "2025-08-25 13:06:00 UTC"
$ws.Cells.Item($newRow, 2).Value = "2025-08-25 18:36:00 IST"
$ws.Cells.Item($newRow, 3).Value = "SKIPPED"
$ws.Cells.Item($newRow, 4).Value = "No change in PDF. Skipping download & Excel update."
$ws.Cells.Item($newRow, 5).Value = "https://nalcoindia.com/wp-content/uploads/2019/01/INGOT-21-08-2025.pdf"
$ws.Cells.Item($newRow, 6).Value = ""
$ws.Cells.Item($newRow, 7).Value = 0
$ws.Cells.Item($newRow, 8).Value = ""
